$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before existing row 47, shifting rows 47:69 down to 49:71.
$ws.Rows("47:48").Insert()

# New row 47 - "Andross" / "Primera"
$ws.Range("A47").Value = 1
$ws.Range("B47").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C47").Value = "Arica y Parinacota"
$ws.Range("D47").Value = 44960
$ws.Range("E47").Value = 15
$ws.Range("F47").Value = "Fruta"
$ws.Range("G47").Value = 100103
$ws.Range("H47").Value = "Frutos de hueso (carozo)"
$ws.Range("I47").Value = 100103004
$ws.Range("J47").Value = "Durazno"
$ws.Range("K47").Value = "Andross"
$ws.Range("L47").Value = "Primera"
$ws.Range("M47").Value = 300
$ws.Range("N47").Value = 25000
$ws.Range("O47").Value = 26000
$ws.Range("P47").Value = 25500
$ws.Range("Q47").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R47").Value = "Región de O'Higgins"
$ws.Range("S47").Value = 1417
$ws.Range("T47").Value = 18

# New row 48 - "Elegant Lady" / "Primera"
$ws.Range("A48").Value = 1
$ws.Range("B48").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C48").Value = "Arica y Parinacota"
$ws.Range("D48").Value = 44960
$ws.Range("E48").Value = 15
$ws.Range("F48").Value = "Fruta"
$ws.Range("G48").Value = 100103
$ws.Range("H48").Value = "Frutos de hueso (carozo)"
$ws.Range("I48").Value = 100103004
$ws.Range("J48").Value = "Durazno"
$ws.Range("K48").Value = "Elegant Lady"
$ws.Range("L48").Value = "Primera"
$ws.Range("M48").Value = 270
$ws.Range("N48").Value = 25000
$ws.Range("O48").Value = 26000
$ws.Range("P48").Value = 25500
$ws.Range("Q48").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R48").Value = "Región de O'Higgins"
$ws.Range("S48").Value = 1417
$ws.Range("T48").Value = 18
